$wb = $excel.ActiveWorkbook

# Colour used by the workbook's existing "HyperLink" font (RGB 6495ED encoded
# the way Font.Color expects it: R + G*256 + B*65536)
$hlColor = 15570276

function Style-AsHyperlink($rng) {
    $rng.Font.Color = $hlColor
    $rng.Font.Underline = 2
}

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

# Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("B3").Value = "Handed back: in sync with en-US"

# New "Latest Target File" (E) / "Latest Handback File" (F) hyperlinked cells
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/b1031b44fae42592383c4043758b4634e62ee0ea/e2e/ef7673a2-0a5b-43bf-bf1d-874a5f3e81b4.md", "", "", "ef7673a2-0a5b-43bf-bf1d-874a5f3e81b4.md")
Style-AsHyperlink $ws.Range("E2")

$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d061033bb1a92308aebed595b54fb33824e6dc90/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/ef7673a2-0a5b-43bf-bf1d-874a5f3e81b4.3a45e9c93f464e6a8fb0a1d1754b1a0308354559.zh-cn.xlf", "", "", "ef7673a2-0a5b-43bf-bf1d-874a5f3e81b4.3a45e9c93f464e6a8fb0a1d1754b1a0308354559.zh-cn.xlf")
Style-AsHyperlink $ws.Range("F2")

$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/b1031b44fae42592383c4043758b4634e62ee0ea/e2e/ef7673a2-0a5b-43bf-bf1d-874a5f3e81b4.md", "", "", "ef7673a2-0a5b-43bf-bf1d-874a5f3e81b4.md")
Style-AsHyperlink $ws.Range("E3")

$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d061033bb1a92308aebed595b54fb33824e6dc90/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/ef7673a2-0a5b-43bf-bf1d-874a5f3e81b4.3a45e9c93f464e6a8fb0a1d1754b1a0308354559.zh-cn.xlf", "", "", "ef7673a2-0a5b-43bf-bf1d-874a5f3e81b4.3a45e9c93f464e6a8fb0a1d1754b1a0308354559.zh-cn.xlf")
Style-AsHyperlink $ws.Range("F3")

# Latest Handback DateTime (G) now has a real timestamp instead of the
# "0001-01-01 00:00:00" placeholder
$ws.Range("G2").Value = "2016-03-02 10:25:56"
$ws.Range("G3").Value = "2016-03-02 10:25:56"

# Handoff Reason (H): "Ignored" -> "Include"
$ws.Range("H2").Value = "Include"
$ws.Range("H3").Value = "Include"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("de-de")

$ws2.Range("B2").Value = "Handed back: in sync with en-US"
$ws2.Range("B3").Value = "Handed back: in sync with en-US"

$ws2.Hyperlinks.Add($ws2.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/b1031b44fae42592383c4043758b4634e62ee0ea/e2e/ef7673a2-0a5b-43bf-bf1d-874a5f3e81b4.md", "", "", "ef7673a2-0a5b-43bf-bf1d-874a5f3e81b4.md")
Style-AsHyperlink $ws2.Range("E2")

$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f5fa7b0b7c01db2fd06cf4452367f24892fac2c5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/ef7673a2-0a5b-43bf-bf1d-874a5f3e81b4.3a45e9c93f464e6a8fb0a1d1754b1a0308354559.de-de.xlf", "", "", "ef7673a2-0a5b-43bf-bf1d-874a5f3e81b4.3a45e9c93f464e6a8fb0a1d1754b1a0308354559.de-de.xlf")
Style-AsHyperlink $ws2.Range("F2")

$ws2.Hyperlinks.Add($ws2.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/b1031b44fae42592383c4043758b4634e62ee0ea/e2e/ef7673a2-0a5b-43bf-bf1d-874a5f3e81b4.md", "", "", "ef7673a2-0a5b-43bf-bf1d-874a5f3e81b4.md")
Style-AsHyperlink $ws2.Range("E3")

$ws2.Hyperlinks.Add($ws2.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f5fa7b0b7c01db2fd06cf4452367f24892fac2c5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/ef7673a2-0a5b-43bf-bf1d-874a5f3e81b4.3a45e9c93f464e6a8fb0a1d1754b1a0308354559.de-de.xlf", "", "", "ef7673a2-0a5b-43bf-bf1d-874a5f3e81b4.3a45e9c93f464e6a8fb0a1d1754b1a0308354559.de-de.xlf")
Style-AsHyperlink $ws2.Range("F3")

$ws2.Range("G2").Value = "2016-03-02 10:26:16"
$ws2.Range("G3").Value = "2016-03-02 10:26:16"

$ws2.Range("H2").Value = "Include"
$ws2.Range("H3").Value = "Include"
